# Reindex CO2 budget fixed.
# The "lifetime"/"years" sub-header (column K) is renamed to "base_level"
# and a brand-new "lifetime"/"years" column is inserted at L, taking over
# the old K-column lifetime values (a few of which were also revised).
# A handful of CAPEX (column C) figures were corrected too, and those
# revised cells get a highlight fill.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("values")

# --- header row 1 & 2: rename K column, introduce L column ---------------
$ws.Cells.Item(1, 11).Value = "base_level"   # K1: lifetime -> base_level
$ws.Cells.Item(1, 12).Value = "lifetime"     # L1: new "lifetime" header

$ws.Cells.Item(2, 11).Clear()                # K2: drop old "years" label
$ws.Cells.Item(2, 12).Value = "years"        # L2: new "years" sub-label

# --- data rows: move the lifetime figures from K to L --------------------

# Row 3 - NGPP/CCGT
$ws.Cells.Item(3, 3).Value = 800             # CAPEX 1100 -> 800
$ws.Cells.Item(3, 11).Clear()
$ws.Cells.Item(3, 12).Value = 30             # lifetime 25 -> 30

# Row 4 - NGPP/OCGT
$ws.Cells.Item(4, 11).Clear()
$ws.Cells.Item(4, 12).Value = 30             # lifetime 25 -> 30

# Row 5 - Nuclear/Uranium
$ws.Cells.Item(5, 11).Value = 0.8            # new base_level value
$ws.Cells.Item(5, 12).Value = 40             # lifetime unchanged (moved)

# Row 6 - Hydro/Reservoir
$ws.Cells.Item(6, 11).Clear()
$ws.Cells.Item(6, 12).Value = 50             # lifetime unchanged (moved)

# Row 7 - Hydro/Run-of-river
$ws.Cells.Item(7, 11).Clear()
$ws.Cells.Item(7, 12).Value = 50             # lifetime unchanged (moved)

# Row 8 - Storage/Pumped-hydro
$ws.Cells.Item(8, 11).Clear()
$ws.Cells.Item(8, 12).Value = 50             # lifetime unchanged (moved)

# Row 9 - Storage/Li-ion
$ws.Cells.Item(9, 3).Value = 422             # CAPEX 300 -> 422
$ws.Cells.Item(9, 4).Value = 36              # FOM 0.54 -> 36
$ws.Cells.Item(9, 5).Value = 0.0027          # VOM 1.7E-3 -> 2.7E-3
$ws.Cells.Item(9, 11).Clear()
$ws.Cells.Item(9, 12).Value = 10             # lifetime unchanged (moved)

# Row 10 - Wind/Onshore
$ws.Cells.Item(10, 3).Value = 1100           # CAPEX 980 -> 1100
$ws.Cells.Item(10, 11).Clear()
$ws.Cells.Item(10, 12).Value = 25            # lifetime 30 -> 25

# Row 11 - Wind/Offshore
$ws.Cells.Item(11, 3).Value = 2200           # CAPEX 1810 -> 2200
$ws.Cells.Item(11, 11).Clear()
$ws.Cells.Item(11, 12).Value = 25            # lifetime 30 -> 25

# Row 12 - Wind/Floating
$ws.Cells.Item(12, 3).Value = 3500           # CAPEX 2263 -> 3500
$ws.Cells.Item(12, 4).Value = 41.5           # FOM unchanged, but re-written (style highlight)
$ws.Cells.Item(12, 5).Value = 0.0025         # VOM unchanged, but re-written (style highlight)
$ws.Cells.Item(12, 11).Clear()
$ws.Cells.Item(12, 12).Value = 25            # lifetime 30 -> 25

# Row 13 - PV/Utility
$ws.Cells.Item(13, 3).Value = 500            # CAPEX 330 -> 500
$ws.Cells.Item(13, 11).Clear()
$ws.Cells.Item(13, 12).Value = 25            # lifetime 40 -> 25

# Row 14 - PV/Residential
$ws.Cells.Item(14, 3).Value = 800            # CAPEX 490 -> 800
$ws.Cells.Item(14, 11).Clear()
$ws.Cells.Item(14, 12).Value = 25            # lifetime 40 -> 25

# Row 15 - Transmission/HVAC_OHL
$ws.Cells.Item(15, 11).Clear()
$ws.Cells.Item(15, 12).Value = 40            # lifetime unchanged (moved)

# Row 16 - Transmission/HVAC_UC
$ws.Cells.Item(16, 11).Clear()
$ws.Cells.Item(16, 12).Value = 40            # lifetime unchanged (moved)

# Row 17 - Transmission/HVDC_OHL
$ws.Cells.Item(17, 11).Clear()
$ws.Cells.Item(17, 12).Value = 40            # lifetime unchanged (moved)

# Row 18 - Transmission/HVDC_UC
$ws.Cells.Item(18, 11).Clear()
$ws.Cells.Item(18, 12).Value = 40            # lifetime unchanged (moved)

# Row 19 - Transmission/HVDC_SC
$ws.Cells.Item(19, 11).Clear()
$ws.Cells.Item(19, 12).Value = 40            # lifetime unchanged (moved)

# Row 20 - Biomass/Wood Residue
$ws.Cells.Item(20, 11).Clear()
$ws.Cells.Item(20, 12).Value = 30            # lifetime 25 -> 30

# Row 21 - Biomass/Municipal Waste
$ws.Cells.Item(21, 11).Clear()
$ws.Cells.Item(21, 12).Value = 30            # lifetime 25 -> 30

# --- highlight the revised CAPEX/FOM/VOM cells with a single new fill ----
# (applied as two contiguous ranges - the COM bridge here only paints the
#  first area of a multi-area/union Range, so issue one call per block)
$ws.Range("C3").Interior.ThemeColor = 4
$ws.Range("C12:E12").Interior.ThemeColor = 4

# --- selection moved to C16 as last user action ---------------------------
$ws.Range("C16").Select()
